$d = $word.ActiveDocument
$wns = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

function Set-ParagraphText {
    param($Index, $NewText, $RunPropsXml)
    $p = $d.Paragraphs.Item($Index)
    $start = $p.Range.Start
    $end = $p.Range.End - 1
    $r = $d.Range($start, $end)
    $escaped = $NewText -replace "&", "&amp;" -replace "<", "&lt;" -replace ">", "&gt;"
    $xml = "<w:p $wns><w:r>$RunPropsXml<w:t xml:space='preserve'>$escaped</w:t></w:r></w:p>"
    $r.InsertXML($xml)
}

# 1. Title heading (Heading 1) -- paragraph 1
Set-ParagraphText 1 "Play Mystery of Eldorado Free Slot - Review" ""

# "What we like" bullet list (paragraphs 36-39)
Set-ParagraphText 36 "Engaging storyline" ""
Set-ParagraphText 37 "Pleasing design" ""
Set-ParagraphText 38 "Exciting special features" ""
Set-ParagraphText 39 "Accessible betting options" ""

# "What we don't like" bullet list (paragraphs 41-42)
Set-ParagraphText 41 "Frequently used theme" ""
Set-ParagraphText 42 "Limited number of pay lines" ""

# Bold restatement of title near the end (paragraph 43)
Set-ParagraphText 43 "Play Mystery of Eldorado Free Slot - Review" "<w:rPr><w:b/></w:rPr>"

# Italic meta description (paragraph 44)
Set-ParagraphText 44 "Play Mystery of Eldorado for free with this review. Discover Mayan treasures and enjoy engaging gameplay." "<w:rPr><w:i/></w:rPr>"
